# Update countries & provincias Spain
# - Re-sort three pairs/groups of countries whose case totals overtook their
#   neighbours in the ranking (Catar/Marruecos, Portugal/Costa Rica,
#   Birmania/Senegal/Zambia) by swapping the country-name cells and the
#   numeric stats that go with them.
# - Refresh the numeric COVID stats (columns B-H) for the affected rows.
# - Bump the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Octubre de 2020 a las 16:52"

# --- Country name swaps (A column) ------------------------------------
$ws.Range("A34").Value = "Catar"
$ws.Range("A35").Value = "Marruecos"

$ws.Range("A52").Value = "Portugal"
$ws.Range("A53").Value = "Costa Rica"

$ws.Range("A91").Value = "Birmania"
$ws.Range("A92").Value = "Senegal"
$ws.Range("A93").Value = "Zambia"

# --- Numeric stats refresh (B:H columns) ------------------------------
# row -> @{ B=...; C=...; D=...; E=...; F=...; G=...; H=... }
$rowData = @{
    4   = @{ B=7501186; C=6515; D=4737664; E=2550792; F=0; G=70; H=212730 }
    15  = @{ B=466590;  C=1840; D=439607;  E=14116;   F=0; G=45; H=12867 }
    18  = @{ B=372259;  C=4785; D=299804;  E=63157;   F=0; G=67; H=9298 }
    25  = @{ B=296615;  C=1085; D=259500;  E=27527;   F=0; G=2;  H=9588 }
    27  = @{ B=258920;  C=5430; D=185122;  E=72165;   F=0; G=11; H=1633 }
    34  = @{ B=126164;  C=205;  D=123108;  E=2841;    F=0; G=1;  H=215 }
    35  = @{ B=126044;  C=0;    D=104136;  E=19679;   F=0; G=0;  H=2229 }
    46  = @{ B=93090;   C=681;  D=81466;   E=8357;    F=0; G=6;  H=3267 }
    52  = @{ B=77284;   C=888;  D=49359;   E=25942;   F=0; G=6;  H=1983 }
    53  = @{ B=76828;   C=0;    D=39843;   E=36068;   F=0; G=0;  H=917 }
    73  = @{ B=38923;   C=210;  D=25114;   E=13084;   F=0; G=7;  H=725 }
    75  = @{ B=35717;   C=509;  D=20334;   E=14813;   F=0; G=11; H=570 }
    76  = @{ B=33735;   C=73;   D=31536;   E=1448;    F=0; G=1;  H=751 }
    80  = @{ B=27975;   C=226;  D=21108;   E=5996;    F=0; G=10; H=871 }
    91  = @{ B=15525;   C=1142; D=4378;    E=10794;   F=0; G=32; H=353 }
    92  = @{ B=15051;   C=32;   D=12694;   E=2045;    F=0; G=1;  H=312 }
    93  = @{ B=14802;   C=0;    D=13961;   E=508;     F=0; G=0;  H=333 }
    97  = @{ B=13650;   C=549;  D=4752;    E=8819;    F=0; G=10; H=79 }
    116 = @{ B=6704;    C=149;  D=2112;    E=4476;    F=0; G=5;  H=116 }
    133 = @{ B=4613;    C=43;   D=2729;    E=1808;    F=0; G=0;  H=76 }
    135 = @{ B=4220;    C=97;   D=2290;    E=1834;    F=0; G=0;  H=96 }
    143 = @{ B=3388;    C=6;    D=3245;    E=130;     F=0; G=0;  H=13 }
}

foreach ($row in $rowData.Keys) {
    $cols = $rowData[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
